# gcs_daten.xlsx edit script
# - Insert a new "zeit" (time) column between "datum" and "gcs_auge"
# - Fill it with half-hour-stepped, hour-incrementing time values (07:30 .. 26:30)
# - Retype a couple of legacy numeric cells in the gcs_* columns as text,
#   matching the rest of the column (data-cleanup / validation pass)
# - Append three new patient observation rows (19-21)
# - Misc cosmetic bits: selection, page setup

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert new column C ("zeit"). Old C/D/E (gcs_auge/verbal/motorisch)
#    shift right to D/E/F.
# ---------------------------------------------------------------------------
$ws.Columns("C").Insert()

# Header cell + column formatting (bold, centered, custom time format)
$ws.Range("C1").Value2 = "zeit"
$ws.Range("C1").Font.Bold = $true
$ws.Range("C1").HorizontalAlignment = -4108
$ws.Range("C1").NumberFormat = "h:mm;@"

# Column width to match column B, and default (non-header) time format
$ws.Columns("C").ColumnWidth = 19.83
$ws.Range("C2:C21").NumberFormat = "h:mm;@"

# ---------------------------------------------------------------------------
# 2. Fill the time column: 07:30, 08:30, ... incrementing by exactly one
#    hour per row, for the 20 data rows (rows 2-21).
# ---------------------------------------------------------------------------
$startTime = 0.3125
for ($i = 0; $i -lt 20; $i++) {
    $fillRow = $i + 2
    $ws.Cells.Item($fillRow, 3).Value2 = $startTime + $i * (1 / 24)
}

# ---------------------------------------------------------------------------
# 3. Data cleanup: rows 17-18 had D/E/F (previously C/D/E) stored as plain
#    numbers while every other row stores them as text. Normalize them to
#    text so the whole column is consistent.
# ---------------------------------------------------------------------------
function Set-TextValue($cell, $text) {
    $targetRange = $ws.Range($cell)
    $targetRange.NumberFormat = "@"
    $targetRange.Value2 = $text
    $targetRange.ClearFormats()
}

Set-TextValue "D17" "3"
Set-TextValue "E17" "5"
Set-TextValue "F17" "5"

Set-TextValue "D18" "4"
Set-TextValue "E18" "5"
Set-TextValue "F18" "5"

# ---------------------------------------------------------------------------
# 4. Append three new observation rows.
# ---------------------------------------------------------------------------
$newRows = @(
    @{ Row = 19; Pid = "4747573"; Datum = 43648; Auge = "4"; Verbal = "5"; Motorisch = "6" },
    @{ Row = 20; Pid = "28374467"; Datum = 43668; Auge = "4"; Verbal = "5"; Motorisch = "6" },
    @{ Row = 21; Pid = "666666";  Datum = 43669; Auge = "3"; Verbal = "5"; Motorisch = "5" }
)

foreach ($entry in $newRows) {
    $rowNum = $entry.Row
    $stepIdx = $rowNum - 2

    Set-TextValue "A$rowNum" $entry.Pid

    $bcell = $ws.Range("B$rowNum")
    $bcell.Value2 = $entry.Datum
    $bcell.NumberFormat = "yyyy\-mm\-dd"

    $ccell = $ws.Range("C$rowNum")
    $ccell.Value2 = $startTime + $stepIdx * (1 / 24)
    $ccell.NumberFormat = "h:mm;@"

    Set-TextValue "D$rowNum" $entry.Auge
    Set-TextValue "E$rowNum" $entry.Verbal
    Set-TextValue "F$rowNum" $entry.Motorisch
}

# ---------------------------------------------------------------------------
# 5. Misc cosmetic updates
# ---------------------------------------------------------------------------
$ws.Range("C23").Select()

$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1
